$d = $word.ActiveDocument

# Locate the trailing paragraph that holds the (hidden) _GoBack bookmark -
# that is the paragraph the new content must be inserted in front of,
# while keeping the bookmark on the final paragraph of the block.
$bm = $d.Bookmarks.Item("_GoBack")
$target = $bm.Range.Paragraphs.Item(1).Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Introducir en el email el código de invitación y la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de la página de verificación,</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Cuando se crea el usuario se tiene que autogenerar el código y guardarlo en el atributo.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Una vez el usuario está en la página de invitación introduce el código y la contraseña.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D97608" w:rsidRDefault="00D97608" w:rsidP="00A94CBE"><w:r><w:t xml:space="preserve">El controlador pedirá al servicio primero que valide el código, y si es válido (de un usuario existente y deshabilitado), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>seteará</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> la contraseña y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>redireccionará</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a la página de login.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

# InsertXML replaces the contents of the target range in place (it does not
# splice in new paragraphs elsewhere), so calling it on the whole bookmark
# paragraph's range swaps that single paragraph for the five paragraphs
# above - an empty paragraph, three new instructional paragraphs, and the
# final paragraph (carrying the original bookmark) that picks up the last
# chunk of new text.
$null = $target.InsertXML($xml)

Write-Output "Inserted invitation/validation paragraphs before the _GoBack bookmark."
